$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new rows needed (top to bottom), shifting existing blocks down ---
# Blank separator row originally at row 13 -> push down, new content row takes its place
$ws.Rows("13").Insert()
# Blank separator row originally at row 16 is now at row 17 after the previous insert
$ws.Rows("17").Insert()
# Blank separator row originally at row 19 is now at row 21 after the two previous inserts
$ws.Rows("21").Insert()

# --- Write the new shared-string values in the exact order they first appear ---
# (this keeps the shared-strings table ordering aligned with the authored edit)
$ws.Range("C27").Value = "LOO Runs (within subject)"
$ws.Range("A28").Value = "slurm-42682826"
$ws.Range("A4").Value = "slurm-42682904"
$ws.Range("A13").Value = "slurm-42682907"
$ws.Range("A17").Value = "slurm-42682908"
$ws.Range("A21").Value = "slurm-42682914"
$ws.Range("A25").Value = "slurm-42682920"

# --- Fill in the rest of the new "LOO Runs (within subject)" block (rows 27-28) ---
$ws.Range("A27").Value = "4 word blocks, 4 word separated"
$ws.Range("B27").Value = "word-level fMRIs averaged"
$ws.Range("D27").Value = 10

# Update the active selection to match the authored workbook state
$ws.Range("B25").Select()
